# Apply the timesheet update for week commencing row 10 (date 43178):
# - Tuesday hours (C10) recorded as 1.5 (was 0)
# - The shared formula in I10 (SUM(B10:H10)) will recalc to 9.75
# - The total in I19 (SUM(I2:I18)) will recalc to 278
# - Selection moves to J20 (cell after the final data row)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Tuesday hours for the week of 43178 (row 10)
$ws.Range("C10").Value = 1.5

# Recalculate so dependent formulas (I10, I19) pick up the new total
$excel.Calculate()

# Move the active selection to J20, matching the saved workbook view
$ws.Range("J20").Select()
